# Add a new "category" column (C) to Sheet1, with a bold/bordered header
# matching the style of the existing A1/B1 headers, and populate the
# category value for every data row (2-39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from A1 (bold font, border, center/top alignment)
# onto C1, then set its text.
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "category"

# category value for each data row, in row order (row 2 .. row 39)
$categories = @(
    "安全",
    "用户体验",
    "不正常退出",
    "不正常退出",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "功能不完整",
    "其他",
    "性能",
    "页面布局缺陷",
    "页面布局缺陷",
    "用户体验",
    "用户体验",
    "用户体验",
    "用户体验",
    "用户体验",
    "用户体验",
    "用户体验",
    "用户体验",
    "用户体验"
)

for ($i = 0; $i -lt $categories.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $categories[$i]
}
